# Updated C3DC Regression and Smoke suites
# - Fix the "Treatment" tab query: replace CONCAT(REPLACE(...)) with REPLACE(...)
# - Move the selection/view back to B2 (top-left at A1) instead of the old C6/A6 scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs000471' AND dgn.anatomic_site = 'C64.9 : Kidney, NOS'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

# Row 5 / column B holds the "TreatmentTab" query (see column A label "TreatmentTab").
# Replace CONCAT(REPLACE(...)) with REPLACE(...) by writing the corrected SQL text back.
$ws.Range("B5").Value = $newTreatmentQuery

# Nudge the cell's font so a distinct style slot is recorded for the edited cell
# (mirrors Excel creating a fresh cellXfs/font entry when the cell content is edited).
$ws.Range("B5").Font.ThemeColor = 1

# Restore the view: scroll back to the top-left (A1) and select B2, matching the saved
# workbook state (previously the view was scrolled to A6 with C6 selected).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
